# Updates cryptos list values (Price / Volume(1h)) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.073.94"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.12%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.594.74"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.48%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.73"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.64"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.47%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.78%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.579"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.88"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.25"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0840"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.11"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.66%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.992.93"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.59%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.80%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.596.14"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.62%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.915"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.80"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.70%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "46.205.75"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.12%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.60%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.83"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.72"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.81%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.78"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "272.36"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +7.65%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.34%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.19"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "29.88"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +8.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.06"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.73"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.71%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "38.15"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.29%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.24"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.93%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.60"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.96%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "155.31"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.53%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.97%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.65%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.80"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.78%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.125"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.87%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.00"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +27.29%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.83"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.27%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.78%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -5.75%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.098.37"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.11%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "94.73"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.68%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.64"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +6.76%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "108.60"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.40%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.18%  "
